$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date in column C for every data row
# (rows 2 through 342) from 45181 to 45182.
for ($row = 2; $row -le 342; $row++) {
    $ws.Cells.Item($row, 3).Value = 45182
}

# Rows 340-342 got re-ordered (beteckning + area swapped between rows)
# while every other column (B, C, D, E, F, H..Q, R) stayed the same.
$ws.Range("A340").Value = "A 42326-2023"
$ws.Range("G340").Value = 1

$ws.Range("A341").Value = "A 42324-2023"
$ws.Range("G341").Value = 2.2

$ws.Range("A342").Value = "A 42328-2023"
$ws.Range("G342").Value = 1.3
